# This script applies a cyclic permutation of the data rows 2-12 on the
# active worksheet. Rows 4 and 13 are unaffected (they map to themselves).
# Only the columns that actually differ between rows are touched:
# A, B, D, E, F, G, H, Q, R, Z, AB, AC
# (columns C, I, K, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY
#  are identical across every affected row, so they are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row r(dest) receives the former content of row r(src).
$rowMap = @{
    2  = 3
    3  = 9
    5  = 12
    6  = 5
    7  = 10
    8  = 6
    9  = 7
    10 = 2
    11 = 8
    12 = 11
}

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

# 1) Snapshot the "before" state of every source row for the columns that
#    vary, plus whether the "AC" (Publik kommentar) cell is populated.
$snapshot = @{}
foreach ($src in ($rowMap.Values | Sort-Object -Unique)) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$src").Value2
    }
    $acCell = $ws.Range("AC$src")
    $acText = $acCell.Value2
    if ($acText -eq $null -or $acText -eq "") {
        $rowValues["AC"] = $null
    } else {
        $rowValues["AC"] = $acText
    }
    $snapshot[$src] = $rowValues
}

# 2) Write the captured values into their destination rows.
foreach ($dest in $rowMap.Keys) {
    $src = $rowMap[$dest]
    $rowValues = $snapshot[$src]

    foreach ($col in $cols) {
        $ws.Range("$col$dest").Value2 = $rowValues[$col]
    }

    if ($rowValues["AC"] -ne $null) {
        $ws.Range("AC$dest").Value2 = $rowValues["AC"]
    } else {
        $ws.Range("AC$dest").ClearContents() | Out-Null
    }
}
